$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("E2").Value = 2
$ws.Range("F2").Value = 1
$ws.Range("G2").Value = 0.5854969999999999
$ws.Range("H2").Value = 1.170994
$ws.Range("M2").Value = 2.8377025
$ws.Range("N2").Value = 5.675405
$ws.Range("O2").Value = 0.08520139853031897
$ws.Range("P2").Value = 0.05979321158534227
$ws.Range("Q2").Value = 1.6614663006425
$ws.Range("R2").Value = 6.645865202569999
$ws.Range("S2").Value = 0.08520139853031897
$ws.Range("T2").Value = 0.05979321158534227
$ws.Range("E3").Value = 2
$ws.Range("F3").Value = 1
$ws.Range("G3").Value = 0.5854969999999999
$ws.Range("H3").Value = 1.170994
$ws.Range("O3").Value = 0.1519846751111432
$ws.Range("P3").Value = 0.1599912441592654
$ws.Range("Q3").Value = 2.963770786243666
$ws.Range("R3").Value = 17.782624717462
$ws.Range("S3").Value = 0.1519846751111432
$ws.Range("T3").Value = 0.1599912441592654
$ws.Range("E4").Value = 2
$ws.Range("F4").Value = 1
$ws.Range("G4").Value = 0.5854969999999999
$ws.Range("H4").Value = 1.170994
$ws.Range("M4").Value = 4.372280666666666
$ws.Range("N4").Value = 13.116842
$ws.Range("O4").Value = 0.1312767732230829
$ws.Range("P4").Value = 0.1381924477702479
$ws.Range("Q4").Value = 2.559957213491332
$ws.Range("R4").Value = 15.359743280948
$ws.Range("S4").Value = 0.1312767732230829
$ws.Range("T4").Value = 0.1381924477702479
$ws.Range("E5").Value = 2
$ws.Range("F5").Value = 1
$ws.Range("G5").Value = 0.5854969999999999
$ws.Range("H5").Value = 1.170994
$ws.Range("M5").Value = 2.1625465
$ws.Range("N5").Value = 4.325093
$ws.Range("O5").Value = 0.06492998691259792
$ws.Range("P5").Value = 0.04556700374251401
$ws.Range("Q5").Value = 1.2661644881105
$ws.Range("R5").Value = 5.064657952441999
$ws.Range("S5").Value = 0.06492998691259792
$ws.Range("T5").Value = 0.04556700374251401
$ws.Range("E6").Value = 2
$ws.Range("F6").Value = 1
$ws.Range("G6").Value = 0.5854969999999999
$ws.Range("H6").Value = 1.170994
$ws.Range("M6").Value = 15.97657333333333
$ws.Range("N6").Value = 47.92972
$ws.Range("O6").Value = 0.479693129114909
$ws.Range("P6").Value = 0.5049634148023289
$ws.Range("Q6").Value = 9.354235756946666
$ws.Range("R6").Value = 56.12541454167999
$ws.Range("S6").Value = 0.479693129114909
$ws.Range("T6").Value = 0.5049634148023289
$ws.Range("E7").Value = 2
$ws.Range("F7").Value = 1
$ws.Range("G7").Value = 0.5854969999999999
$ws.Range("H7").Value = 1.170994
$ws.Range("M7").Value = 2.894743333333333
$ws.Range("N7").Value = 8.68423
$ws.Range("O7").Value = 0.08691403710794816
$ws.Range("P7").Value = 0.0914926779403015
$ws.Range("Q7").Value = 1.694863537436666
$ws.Range("R7").Value = 10.16918122462
$ws.Range("S7").Value = 0.08691403710794816
$ws.Range("T7").Value = 0.0914926779403015
